# template con colores correctos
# Recolor the "85A0FE" (blue) and "FE938C" (salmon) funnel-step circles to
# "FFD14C" (yellow) on slides 1 and 2, matching the target template.

$p = $ppt.ActivePresentation

# New fill color: FFD14C  (PowerPoint RGB long = B*65536 + G*256 + R)
$newColor = 0x4C * 65536 + 0xD1 * 256 + 0xFF

function Set-FunnelColors {
    param($groupShape, $indices)
    foreach ($idx in $indices) {
        $item = $groupShape.GroupItems.Item($idx)
        $item.Fill.ForeColor.RGB = $newColor
    }
}

# --- Slide 1 ---
$s1 = $p.Slides.Item(1)
Set-FunnelColors $s1.Shapes.Item(2) @(1, 3)   # Group 6
Set-FunnelColors $s1.Shapes.Item(8) @(1, 3)   # Group 3
Set-FunnelColors $s1.Shapes.Item(9) @(1, 3)   # Group 16

# --- Slide 2 ---
$s2 = $p.Slides.Item(2)
Set-FunnelColors $s2.Shapes.Item(2) @(1, 3)   # Group 3
Set-FunnelColors $s2.Shapes.Item(8) @(7, 8)   # Group 4
